$d = $word.ActiveDocument

# --- 1) Cambria -> Ubuntu for the five header-block paragraphs
#        ("Batch: T2", "Assignment No: ", "Title of Assignment: ",
#         "Student Name: ...", "Student PRN: ...") ---
for ($i = 1; $i -le 5; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Name = "Ubuntu"
}

# --- 2) Collapse the "Assignment" / " No" / ": " runs into a single run
#        reading "Assignment No: " ---
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Find.Execute("Assignment No: ", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "Assignment No: ", 2)

# --- 3) Header: split "Walchand College of Engineering, Sangli" into
#        three runs, wrapping the two proper nouns in spell-check
#        proofErr markers, same run formatting throughout ---
$hdr = $d.Sections.Item(1).Headers.Item(1)
$rng = $hdr.Range.Duplicate
$rng.Find.Execute("Walchand College of Engineering, Sangli")

$rPr = '<w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="18"/><w:lang w:val="en-IN"/></w:rPr>'
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
               '<pkg:xmlData>' + `
                   '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
                       '<w:body>' + `
                           '<w:p>' + `
                               '<w:proofErr w:type="spellStart"/>' + `
                               '<w:r>' + $rPr + '<w:t>Walchand</w:t></w:r>' + `
                               '<w:proofErr w:type="spellEnd"/>' + `
                               '<w:r>' + $rPr + '<w:t xml:space="preserve"> College of Engineering, </w:t></w:r>' + `
                               '<w:proofErr w:type="spellStart"/>' + `
                               '<w:r>' + $rPr + '<w:t>Sangli</w:t></w:r>' + `
                               '<w:proofErr w:type="spellEnd"/>' + `
                           '</w:p>' + `
                       '</w:body>' + `
                   '</w:document>' + `
               '</pkg:xmlData>' + `
           '</pkg:part>' + `
       '</pkg:package>'

$rng.InsertXML($xml)
